$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 23, shifting existing rows 23:108 down to 24:109
$ws.Rows.Item(23).Insert()

# Fill the new row 23 with data
$ws.Cells.Item(23, 1).Value = 6
$ws.Cells.Item(23, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(23, 3).Value = "Metropolitana"
$ws.Cells.Item(23, 4).Value = 44487
$ws.Cells.Item(23, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(23, 5).Value = 13
$ws.Cells.Item(23, 6).Value = 100112001
$ws.Cells.Item(23, 7).Value = "Berenjena"
$ws.Cells.Item(23, 8).Value = "Sin especificar"
$ws.Cells.Item(23, 9).Value = "Primera"
$ws.Cells.Item(23, 10).Value = 2400
$ws.Cells.Item(23, 11).Value = 5000
$ws.Cells.Item(23, 12).Value = 6000
$ws.Cells.Item(23, 13).Value = 5667
$ws.Cells.Item(23, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(23, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(23, 16).Value = 113
$ws.Cells.Item(23, 17).Value = 50
$ws.Cells.Item(23, 18).Value = "Hortaliza"
